$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'36.652.53"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +3.71%  '
$ws.Range('D3').Value = "'1.917.00"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('D5').Value = "'250.40"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.02%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +2.82%  '
$ws.Range('D9').Value = "'0.371"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.94%  '
$ws.Range('D10').Value = "'58.47"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +9.75%  '
$ws.Range('D11').Value = "'0.0765"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.03%  '
$ws.Range('D13').Value = "'14.60"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +8.30%  '
$ws.Range('D14').Value = "'0.817"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +6.81%  '
$ws.Range('D15').Value = "'2.193.32"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('D16').Value = "'5.14"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.13%  '
$ws.Range('D17').Value = "'1.917.86"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.22%  '
$ws.Range('D18').Value = "'36.627.02"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.66%  '
$ws.Range('D19').Value = "'74.63"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.43%  '
$ws.Range('E20').Value = '  +4.97%  '
$ws.Range('D21').Value = "'250.63"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.57%  '
$ws.Range('E22').Value = '  +4.42%  '
$ws.Range('E23').Value = '  +3.55%  '
$ws.Range('D24').Value = "'2.63"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D26').Value = "'2.21"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.23%  '
$ws.Range('D27').Value = "'168.81"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.32%  '
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('E29').Value = '  +2.56%  '
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('D31').Value = "'4.58"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +6.73%  '
$ws.Range('D32').Value = "'0.0621"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.79%  '
$ws.Range('D33').Value = "'4.34"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.95%  '
$ws.Range('D34').Value = "'0.0887"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +22.11%  '
$ws.Range('E35').Value = '  +3.53%  '
$ws.Range('E37').Value = '  +6.64%  '
$ws.Range('E38').Value = '  +3.43%  '
$ws.Range('D39').Value = "'17.87"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +51.01%  '
$ws.Range('E40').Value = '  +4.15%  '
$ws.Range('D41').Value = "'106.50"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +9.99%  '
$ws.Range('E42').Value = '  +4.02%  '
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('D44').Value = "'2.94"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +22.96%  '
$ws.Range('E45').Value = '  +3.85%  '
$ws.Range('D46').Value = "'1.346.75"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.09%  '
$ws.Range('D47').Value = "'2.39"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.86%  '
$ws.Range('D48').Value = "'0.0816"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('D49').Value = "'2.81"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.80%  '
$ws.Range('D50').Value = "'43.75"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.26%  '
$ws.Range('E51').Value = '  +1.64%  '
